$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of author data (row 3)
$ws.Range("A3").Value = "prueba"
$ws.Range("B3").Value = "final"
$ws.Range("C3").Value = "F"
$ws.Range("D3").Value = "V585"
$ws.Range("E3").Value = "a@gmail.com"
$ws.Range("F3").Value = 54454545455
$ws.Range("G3").Value = 45454545455
$ws.Range("H3").Value = "CX"
$ws.Range("I3").Value = "N"
$ws.Range("J3").Value = "Guarenas"
$ws.Range("K3").Value = "Licenciado"
$ws.Range("L3").Value = "Ok"
$ws.Range("M3").Value = "humanidades"
$ws.Range("N3").Value = "derecho"
$ws.Range("Q3").Value = "universidad de prueba"

# Style the email cell like a hyperlink (blue font) and attach the mailto hyperlink
$ws.Range("E3").Font.Color = 16711680
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:a@gmail.com", "", "", "a@gmail.com")

# Move selection to the row below the newly inserted data
$ws.Range("A4:Q4").Select()
